$wb = $excel.ActiveWorkbook

# Update data values on the BVTStL sheet (sheet2): vehicle types are no
# longer marked as subject to LCFS (LDVs, HDVs, rail, ships, motorbikes)
$ws = $wb.Worksheets.Item("BVTStL")

$ws.Range("B2:C3").Value = 0
$ws.Range("B5:C7").Value = 0

# Activate the BVTStL sheet (becomes the selected tab / activeTab)
$ws.Activate()

# Update selection on the BVTStL sheet
$ws.Range("C6").Select()
